$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Не"
$ws.Range("F2").Value = "Не"
$ws.Range("E3").Value = "Не"
$ws.Range("E4").Value = "Не"
$ws.Range("F4").Value = "Не"
$ws.Range("E5").Value = "Не"
$ws.Range("F5").Value = "Да"
$ws.Range("E6").Value = "Не"
$ws.Range("F6").Value = "Не"
$ws.Range("L6").Value = "Максимум возраст на клиентот при аплицирање/достасување на кредитот `n20г. на денот на аплицирање`n67г. на денот на достасување ( 70 со вклучен ко-кредитобарател)`n11,60% за кредит до 300.000 МКД`n10.70% за кредит над 300.001 МКД`n1) апликативен трошок 600мкд`n2) административен трошок 2%`n3) полиса за осигурување од незгода 600мкд`nтрошоци за Нотарска солемнизација на Договорот за кредит`n"
$ws.Range("E7").Value = "Не"
$ws.Range("F7").Value = "Не"
$ws.Range("L7").Value = "Максимум возраст на клиентот при аплицирање/достасување на кредитот `n20г. на денот на аплицирање`n67г. на денот на достасување ( 70 со вклучен ко-кредитобарател)`n11,60% за кредит до 300.000 МКД`n10.70% за кредит над 300.001 МКД`n1) апликативен трошок 600мкд`n2) административен трошок 2%`n3) полиса за осигурување од незгода зависно од параметрите на кредитот`nтрошоци за Нотарска солемнизација на Договорот за кредит`n"
$ws.Range("E8").Value = "Не"
$ws.Range("F8").Value = "Не"
$ws.Range("E9").Value = "Не"
$ws.Range("F9").Value = "Не"
$ws.Range("E12").Value = "Не"
$ws.Range("L14").Value = "3,65% фиксна за првите 5 години`nво зависност од обезбедувањето, односно 70% од проценетата вредност на колатералот;`n- потоа 3,7% + 6 месечен еурибор, минимум 6%"
$ws.Range("L15").Value = "3,95 фиксна за првите 10 години `n- потоа 3,7% + 6 месечен еурибор, минимум 6%`nво зависност од обезбедувањето, односно 70% од проценетата вредност на колатералот;"

$ws.Range("A7").Select()
